$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldName = "Y4_B2526_General_&_Special_surgery_1_reference_data_D16092025T112941.xlsx"
$newName = "Y4_B2526_General_&_Special_surgery_1_reference_data_D21092025T123610.xlsx"

for ($r = 2; $r -le 333; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Text -eq $oldName) {
        $cell.Value = $newName
    }
}
